# Helper: convert a "RRGGBB" hex string into the Long value expected by
# PowerPoint's RGBColor.RGB property (R + G*256 + B*65536 — VBA's RGB()
# byte order), since PowerPoint COM stores/returns colors that way.
function Convert-HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on slide 16 ("Google Shape;213;p29") switches to a different
#    built-in table style.
# ---------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{42B9CE19-1C68-40A2-84F9-8228A6210FB7}")
    }
}

# ---------------------------------------------------------------------
# 2) The presentation theme swaps its 12 scheme colors from the
#    "Integral" palette to the standard "Office Theme" palette.
# ---------------------------------------------------------------------
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToVbaRgb $officeThemeColors[$i - 1]
}
